$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-04 Thursday" "2025-12-05 Friday"

Replace-Text "880÷6=" "105÷7="
Replace-Text "847÷8=" "473÷8="
Replace-Text "219÷5=" "424÷6="
Replace-Text "301÷3=" "963÷8="
Replace-Text "283÷9=" "576÷2="

Replace-Text "664÷6=" "237÷2="
Replace-Text "499÷2=" "711÷7="
Replace-Text "917÷6=" "382÷2="
Replace-Text "830÷6=" "462÷5="
Replace-Text "802÷2=" "558÷8="

Replace-Text "418÷5=" "841÷4="
Replace-Text "376÷7=" "574÷8="
Replace-Text "280÷2=" "641÷8="
Replace-Text "906÷5=" "270÷7="
Replace-Text "904÷4=" "207÷4="

Replace-Text "308÷2=" "474÷8="
Replace-Text "737÷9=" "296÷5="
Replace-Text "897÷6=" "322÷4="
Replace-Text "883÷9=" "692÷5="
Replace-Text "313÷3=" "425÷7="

Replace-Text "991÷5=" "610÷7="
Replace-Text "441÷6=" "106÷4="
Replace-Text "747÷9=" "645÷3="
Replace-Text "237÷8=" "965÷9="
Replace-Text "734÷5=" "807÷9="
